$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

# Copy header style from D1 (bold/border) onto the new E1 header cell, then set its text
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("E1").Value = "users"

# Add the "users" list for each project row
$ws.Range("E2").Value = "['Arun Lakshmanan', 'Mitchell Jones']"
$ws.Range("E3").Value = "['Won Dong Shin']"
$ws.Range("E4").Value = "['Mihary Ito']"
$ws.Range("E5").Value = "['Gavin Ananda']"
$ws.Range("E6").Value = "['Aygen Berk Cagilci', 'Daniel Olivas Hernandez']"
$ws.Range("E7").Value = "['Harshal Maske']"
